$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 313
$ws1.Range("F6").Value = 275
$ws1.Range("F7").Value = 147
$ws1.Range("F13").Value = 3246
$ws1.Range("F14").Value = 101
$ws1.Range("F16").Value = 63
$ws1.Range("F17").Value = 30
$ws1.Range("F19").Value = 561
$ws1.Range("F20").Value = 33
$ws1.Range("F21").Value = 645
$ws1.Range("F22").Value = 193
$ws1.Range("F23").Value = 107
$ws1.Range("F25").Value = 40
$ws1.Range("F26").Value = 58
$ws1.Range("F27").Value = 2314
$ws1.Range("F28").Value = 4859
$ws1.Range("F32").Value = 1256
$ws1.Range("F33").Value = 258
$ws1.Range("F34").Value = 2170
$ws1.Range("F36").Value = 479
$ws1.Range("F37").Value = 71
$ws1.Range("F38").Value = 67
$ws1.Range("F39").Value = 149
$ws1.Range("F42").Value = 759
$ws1.Range("F43").Value = 18
$ws1.Range("F44").Value = 447
$ws1.Range("F45").Value = 24

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 314
$ws4.Range("F6").Value = 275
$ws4.Range("F7").Value = 147
$ws4.Range("F13").Value = 3246
$ws4.Range("F14").Value = 102
$ws4.Range("F16").Value = 63
$ws4.Range("F18").Value = 30
$ws4.Range("F20").Value = 561
$ws4.Range("F21").Value = 33
$ws4.Range("F22").Value = 645
$ws4.Range("F23").Value = 193
$ws4.Range("F24").Value = 107
$ws4.Range("F26").Value = 40
$ws4.Range("F27").Value = 58
$ws4.Range("F28").Value = 2314
$ws4.Range("F29").Value = 4859
$ws4.Range("F33").Value = 1256
$ws4.Range("F34").Value = 258
$ws4.Range("F35").Value = 2170
$ws4.Range("F37").Value = 479
$ws4.Range("F38").Value = 71
$ws4.Range("F39").Value = 67
$ws4.Range("F40").Value = 149
$ws4.Range("F43").Value = 759
$ws4.Range("F44").Value = 18
$ws4.Range("F45").Value = 447
$ws4.Range("F46").Value = 24
